$d = $word.ActiveDocument

# --- Change 1: "I could not be prouder..." paragraph ---
# Remove the leading "I could not be prouder to be a part of this group! " clause.
$d.Content.Find.Execute(
    "I could not be prouder to be a part of this group! ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2) | Out-Null

# Tighten the "remote with issues..." clause down to "remote."
$d.Content.Find.Execute(
    "remote with issues such as limitations of technology, scheduling meetings and time zone differences.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "remote.", 2) | Out-Null

# --- Change 2: "In general, our team communicated..." -> "Our team communicated..." ---
$d.Content.Find.Execute(
    "In general, our team communicated", $true, $false, $false, $false, $false,
    $true, 1, $false, "Our team communicated", 2) | Out-Null

# --- Change 3: "We had most of our tasks" -> "We set most of our tasks" ---
$d.Content.Find.Execute(
    "We had most of our tasks", $true, $false, $false, $false, $false,
    $true, 1, $false, "We set most of our tasks", 2) | Out-Null

# --- Change 4: remove " in this instance" ---
$d.Content.Find.Execute(
    "for some in this instance, we should perhaps", $true, $false, $false, $false, $false,
    $true, 1, $false, "for some, we should perhaps", 2) | Out-Null

# --- Change 5: trim the "Everyone was extremely supportive..." paragraph ---
# Drop the "I initially thought..." sentence.
$d.Content.Find.Execute(
    "they could. I initially thought that because we were all complete strangers at the start of the assignment that we would all be cold towards each other and just focus on the task. Everyone was sincere",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "they could. Everyone was sincere", 2) | Out-Null

# Drop the trailing "When we all received..." sentence AND merge the now-empty
# Heading2 paragraph that used to follow it into this paragraph (by including
# the paragraph mark, char 13, at the end of the search text but not the
# replacement text).
$gitHubTail = "about using Git and GitHub. When we all received our results for Assignment 1 everyone praised each other" + [char]8217 + "s results, and everyone was genuinely happy with how we had performed." + [char]13
$d.Content.Find.Execute(
    $gitHubTail, $true, $false, $false, $false, $false,
    $true, 1, $false, "about using Git and GitHub. ", 2) | Out-Null

# --- Change 6: "Morale and a positive work environment..." paragraph ---
$d.Content.Find.Execute(
    "responded positively to each other I feel everyone responded to tasks more willingly.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "responded positively to each other, everyone did their tasks willingly.", 2) | Out-Null
